# LQA_Tester_ProgressTracker.xlsx :: TOTAL sheet
# Simplify the "EN TESTER STATS" table:
#  - drop the "Completion" (old B) and "Actual Issues" (old C) columns
#  - shift the remaining columns (old D..K) left into (new B..I)
#  - rename header "Total" -> "Done"
#  - shrink the title merges from A1:K1 / A14:K14 to A1:I1 / A14:I14
#
# The "EN CATEGORY BREAKDOWN" (rows 18-30) and "EN RANKING" (rows 33-43)
# tables below are untouched by the source diff, so this script only
# touches rows 2-15 of the TOTAL sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TOTAL")

$xlPasteFormats = -4122

$dataRows = @(3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 15)

# Pairs of (destination column, source column) to process strictly left to
# right. Because source = destination + 2 always, and we write destinations
# in increasing column order, every source is read (copy/value) before it
# is ever used as a destination later in the same loop - no temp needed.
$colPairs = @(
    @(2, 4),   # B <- D
    @(3, 5),   # C <- E
    @(4, 6),   # D <- F
    @(5, 7),   # E <- G
    @(6, 8),   # F <- H
    @(7, 9),   # G <- I
    @(8, 10),  # H <- J
    @(9, 11)   # I <- K
)

foreach ($r in $dataRows) {
    foreach ($pair in $colPairs) {
        $destCol = $pair[0]
        $srcCol = $pair[1]

        $srcCell = $ws.Cells.Item($r, $srcCol)
        $destCell = $ws.Cells.Item($r, $destCol)

        # Carry the source's number/fill/border formatting onto the
        # destination (reuses the existing style record instead of
        # fabricating a new one).
        $srcCell.Copy()
        $destCell.PasteSpecial($xlPasteFormats)

        # Carry the value.
        $destCell.Value = $srcCell.Value2
    }

    # The old J/K columns are no longer part of the table.
    $ws.Cells.Item($r, 10).ClearContents()
    $ws.Cells.Item($r, 10).ClearFormats()
    $ws.Cells.Item($r, 11).ClearContents()
    $ws.Cells.Item($r, 11).ClearFormats()
}

$excel.CutCopyMode = $false

# Header row (row 2) labels.
$ws.Cells.Item(2, 2).Value = "Done"
$ws.Cells.Item(2, 3).Value = "Issues"
$ws.Cells.Item(2, 4).Value = "No Issue"
$ws.Cells.Item(2, 5).Value = "Blocked"
$ws.Cells.Item(2, 6).Value = "Fixed"
$ws.Cells.Item(2, 7).Value = "Reported"
$ws.Cells.Item(2, 8).Value = "Checking"
$ws.Cells.Item(2, 9).Value = "Pending"
$ws.Cells.Item(2, 10).ClearContents()
$ws.Cells.Item(2, 10).ClearFormats()
$ws.Cells.Item(2, 11).ClearContents()
$ws.Cells.Item(2, 11).ClearFormats()

# Shrink the section-title merges to match the narrower (9-column) table.
$ws.Range("A1:K1").UnMerge()
$ws.Range("A1:I1").Merge()

$ws.Range("A14:K14").UnMerge()
$ws.Range("A14:I14").Merge()
